$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "{rated_speed} RPM, {pump_hp} BHP {curve_details}" -> insert "pump_"
#    right after the opening brace, so the placeholder becomes
#    "{pump_rated_speed} RPM, {pump_hp} BHP {curve_details}".
#    Word splits the original single run into three runs:
#      "{"  /  "pump_"  /  "rated_speed} RPM, {pump_hp} BHP {curve_details}"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("{rated_speed} RPM, {pump_hp} BHP {curve_details}", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found1) {
    $insPos1 = $rng1.Start + 1
    $ins1 = $d.Range($insPos1, $insPos1)
    $ins1.InsertBefore("pump_")
    $splitRng1 = $d.Range($insPos1, $insPos1 + 5)
    $splitRng1.Bold = 1
    $splitRng1.Bold = 0
}

# ---------------------------------------------------------------------------
# 2) "RATED FOR {rated_hp} HP @ {rated_speed} RPM" -> insert "engine_" right
#    after the second opening brace, so the placeholder becomes
#    "RATED FOR {rated_hp} HP @ {engine_rated_speed} RPM".
#    Word splits the original single run into three runs:
#      "RATED FOR {rated_hp} HP @ {"  /  "engine_"  /  "rated_speed} RPM"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("RATED FOR {rated_hp} HP @ {rated_speed} RPM", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found2) {
    $insPos2 = $rng2.Start + 27
    $ins2 = $d.Range($insPos2, $insPos2)
    $ins2.InsertBefore("engine_")
    $splitRng2 = $d.Range($insPos2, $insPos2 + 7)
    $splitRng2.Bold = 1
    $splitRng2.Bold = 0
}

# ---------------------------------------------------------------------------
# 3) Move the <w:lastRenderedPageBreak/> marker from the row-42 "1" cell to
#    the row-40 "1" cell of the first table (both are Qty cells holding the
#    text "1"). We rebuild each paragraph via InsertXML (which lands just
#    before the existing one) and then delete the now-redundant original.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# 3a) add the marker to row 40
$cellAdd = $t.Cell(40, 1)
$rngAdd = $cellAdd.Range
$collapsedAdd = $d.Range($rngAdd.Start, $rngAdd.Start)
$xmlAdd = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='4946FBF6' w14:textId='3D9A77C5' w:rsidR='00106A97' w:rsidRPr='00106A97' w:rsidRDefault='00106A97' w:rsidP='009E582B'><w:pPr><w:jc w:val='center'/><w:rPr><w:rFonts w:ascii='Montserrat' w:hAnsi='Montserrat'/><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Montserrat' w:hAnsi='Montserrat'/><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr><w:lastRenderedPageBreak/><w:t>1</w:t></w:r></w:p>"
$collapsedAdd.InsertXML($xmlAdd)
$cellAddAfter = $t.Cell(40, 1)
$dupPara1 = $cellAddAfter.Range.Paragraphs.Item(2)
$dupPara1.Range.Delete()

# 3b) remove the marker from row 42
$cellDel = $t.Cell(42, 1)
$rngDel = $cellDel.Range
$collapsedDel = $d.Range($rngDel.Start, $rngDel.Start)
$xmlDel = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='05CF5C22' w14:textId='78231A65' w:rsidR='00106A97' w:rsidRPr='00106A97' w:rsidRDefault='00106A97' w:rsidP='009E582B'><w:pPr><w:jc w:val='center'/><w:rPr><w:rFonts w:ascii='Montserrat' w:hAnsi='Montserrat'/><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Montserrat' w:hAnsi='Montserrat'/><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr><w:t>1</w:t></w:r></w:p>"
$collapsedDel.InsertXML($xmlDel)
$cellDelAfter = $t.Cell(42, 1)
$dupPara2 = $cellDelAfter.Range.Paragraphs.Item(2)
$dupPara2.Range.Delete()

Write-Host "done"
